$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Content edits (shared-string text changes) ---

# "gimme_boxes, waitOnConn" -> "gimme_boxes, waitOnTransport"
$ws.Cells.Item(7, 1).Value = "gimme_boxes, waitOnTransport"

# Rich-text comment cell (F3) loses its bold "How to signify..." run and becomes
# a single plain-text run (also drops the trailing space before the old run).
$ws.Cells.Item(3, 6).Value = "There's a difference between ""client"" and ""server"" because client is: initiating all the connections, generally has more memory/CPU time/TIME_WAIT slots than server, and because the server doesn't really care about all the client's problems with decoding frames. Server is still interested in this ""problems"" data, so it could be sent in a batch."

# gimme_boxes explanation cell (F7) gets an extra trailing sentence.
$ws.Cells.Item(7, 6).Value = "gimme_boxes is sent C2S only to indicate that client wants to start receiving boxes. It is implied that server always wants to receive boxes over any existing or new transport; making the client wait for a ""gimme_boxes"" from the server would introduce a round trip. If waitOnTransport == -1, give client boxes immediately, else, wait for transport #<waitOnTransport> to close first."

# Row 7 no longer has a hard-coded height; Excel re-measures it once the cell
# text got longer (it wraps onto one extra line).
$ws.Rows.Item(7).RowHeight = 38.25

# --- View state ---
$ws.Range("F7").Select() | Out-Null
